# Update UnitMass ("C" column) values in the "+ loading" and "- loading"
# tables on Sheet1 of LoadingPC5.xlsx to reflect the corrected/re-indexed
# DataFrame values described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "C2"  = 39
    "C3"  = 27
    "C4"  = 24
    "C5"  = 28
    "C6"  = 102
    "C7"  = 58
    "C8"  = 23
    "C9"  = 25
    "C10" = 7
    "C11" = 97
    "C12" = 85
    "C13" = 123
    "C14" = 31
    "C15" = 72
    "C16" = 74
    "C17" = 94
    "C18" = 70
    "C19" = 46
    "C20" = 60
    "C21" = 138
    "C23" = 91
    "C24" = 57
    "C25" = 43
    "C26" = 108
    "C27" = 106
    "C28" = 32
    "C29" = 92
    "C30" = 153
    "C31" = 40
    "C32" = 90
    "C33" = 29
    "C34" = 83
    "C35" = 65
    "C36" = 175
    "C37" = 113
    "C38" = 119
    "C39" = 103
    "C40" = 167
    "C41" = 56
    "C42" = 51
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
